# Glenn Maxwell (Kings XI Punjab) innings log — correct the per-innings
# runs / balls / fours figures in columns C, D, E for rows 2-12.
# Values are kept as text (matching the sheet's existing
# "numbers stored as text" convention) by forcing a Text number format
# before writing each cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $value) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# Row 2
Set-TextValue "C2" "0"
Set-TextValue "D2" "2"
Set-TextValue "E2" "0"

# Row 3
Set-TextValue "C3" "32"
Set-TextValue "D3" "24"
Set-TextValue "E3" "3"

# Row 4
Set-TextValue "C4" "7"
Set-TextValue "D4" "12"

# Row 5
Set-TextValue "C5" "10"
Set-TextValue "D5" "5"
Set-TextValue "E5" "2"

# Row 6
Set-TextValue "C6" "6"
Set-TextValue "D6" "6"
Set-TextValue "E6" "1"

# Row 7
Set-TextValue "C7" "11"
Set-TextValue "D7" "7"
Set-TextValue "E7" "1"

# Row 8
Set-TextValue "C8" "1"
Set-TextValue "D8" "4"
Set-TextValue "E8" "0"

# Row 9
Set-TextValue "C9" "11"
Set-TextValue "D9" "18"
Set-TextValue "E9" "0"

# Row 10
Set-TextValue "C10" "13"
Set-TextValue "D10" "9"
Set-TextValue "E10" "2"

# Row 11
Set-TextValue "C11" "12"
Set-TextValue "D11" "13"

# Row 12
Set-TextValue "C12" "5"
Set-TextValue "D12" "6"
